$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being updated per row, in order: B, C, E, F, G, I, L, M, O
$cols = @(2, 3, 5, 6, 7, 9, 12, 13, 15)

$data = @(
    @(15.24594150502986, 10.78029317107008, 12.2716143256284, 16.86991607391233, 3.63545879023024, 22.04260876646949, 10.03568267277051, 14.60136761148097, 21.64217706859149),
    @(14.64847755712248, 10.52385880621348, 12.32233011089698, 15.89584955866808, 3.637418412451142, 22.21166831311537, 10.04526925676741, 14.46538911247506, 21.78665297732572),
    @(14.26930731802055, 10.36249423100027, 12.35513536630354, 15.26997757108491, 3.638684810112848, 22.32100082394035, 10.05257892968026, 14.38270015286681, 21.88164429386773),
    @(14.11189893061301, 10.29582154562277, 12.36892382491229, 15.00819731993403, 3.639216814896257, 22.36694749521824, 10.05591602755954, 14.34923359824358, 21.92193023011123),
    @(14.08559312490176, 10.28469739884723, 12.37123879425823, 14.96433081551593, 3.639306117935939, 22.37466109753307, 10.05649180389803, 14.34369119883178, 21.92871477514984),
    @(14.26719587906251, 10.36159867227624, 12.35531961986807, 15.26647399323137, 3.638691920319268, 22.32161483426796, 10.05262248353134, 14.38224784367458, 21.88218122678813),
    @(15.04260999660249, 10.69272659664646, 12.2887561852346, 16.53996406344768, 3.636121384714259, 22.09975434390255, 10.03869294499518, 14.55433265911385, 21.69068624272112),
    @(16.45755437602974, 11.3082197127296, 12.171384583474, 19.00274580682531, 3.631579608120049, 21.70844854223809, 10.02265234782705, 14.89690895239655, 21.36517751457645),
    @(17.42373765164852, 11.73634208906358, 12.09309472905348, 20.67494806633232, 3.628543787922732, 21.44747988758631, 10.01771230066688, 15.150038226863, 21.15674347140891),
    @(17.84584331795234, 11.92523508841321, 12.05918658860557, 21.3917225636224, 3.627227397081231, 21.33448613860461, 10.01694376529097, 15.26514804122504, 21.06864308785926),
    @(18.00307672434247, 11.99587309895344, 12.04659055510876, 21.65686569030329, 3.626738154679594, 21.29251905268916, 10.01686459825347, 15.30870157347889, 21.03625235903745),
    @(17.96933123108338, 11.98070032806193, 12.04929249126374, 21.60004134736742, 3.626843111282554, 21.30152091797053, 10.01687223894536, 15.29932362313364, 21.04318501603802),
    @(17.85883183530776, 11.93106465462554, 12.05814541673805, 21.4136618050453, 3.627186961767991, 21.3310170323523, 10.01693301089945, 15.26873211140506, 21.06595879400182),
    @(17.79080517188987, 11.90054387536502, 12.06359986378998, 21.29868154950795, 3.627398782895985, 21.34919115045578, 10.01699780157858, 15.24998833847514, 21.08003499189589),
    @(17.39579232340326, 11.72387526812859, 12.09534494495288, 20.62722412089977, 3.628631113528993, 21.45497929687857, 10.01779221719677, 15.14251212644232, 21.16263663517289),
    @(17.1489248819895, 11.61395739634274, 12.11525574922827, 20.20408069597325, 3.629403626167985, 21.52134128717003, 10.01865777238656, 15.07654704086452, 21.21503406875505),
    @(17.00529920193973, 11.5501864243109, 12.12686858937084, 19.95656407809801, 3.62985403997212, 21.56004964230148, 10.01929485211979, 15.03860370143209, 21.24580364164551),
    @(16.95639279725185, 11.52850191047459, 12.13082812742991, 19.87204792380568, 3.630007588849665, 21.5732481992364, 10.01953449255408, 15.02575733016126, 21.25633007815746),
    @(17.17537426597861, 11.62571553943958, 12.11311958931295, 20.24955283636154, 3.629320761440627, 21.51422119948533, 10.01855122760998, 15.08356956111178, 21.20939083149327),
    @(17.89135974812475, 11.94566842228054, 12.05553847812392, 21.46857628470577, 3.627085713968551, 21.32233103427094, 10.01690941745002, 15.2777187981565, 21.05924319059063),
    @(18.34405304513242, 12.14955911658769, 12.01932901694822, 22.22866616901552, 3.625678854356159, 21.20170556971587, 10.0170708692813, 15.40438261973974, 20.96677504442988),
    @(18.10386771377333, 12.04123119496709, 12.0385248400294, 21.82633154458858, 3.626424807488866, 21.26564826572414, 10.01687202243346, 15.33681015268498, 21.01560725997441),
    @(17.16342178122586, 11.62040147872257, 12.11408483050733, 20.22900810905287, 3.629358204999162, 21.51743845917455, 10.0185989620406, 15.08039473697946, 21.21194012843346),
    @(16.08708983127496, 11.14570884051679, 12.20173614876632, 18.34778573295695, 3.63275518613237, 21.80963841441556, 10.02578733829063, 14.80385839751974, 21.44785900968401)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowVals[$j]
    }
}

Write-Output "Updated loading_percent values for rows 2-25"
